# EnvironmentExposureValidation.xlsx -- fix spelling / rebrand "BioGears" -> "Engine"/"Pulse"
#
# This script applies the content edits described by the commit:
#   "Fixed spelling errors in validation tables."
# which, on closer inspection of the OOXML diff, is really a rebrand of the
# "BioGears" engine name to "Engine"/"Pulse" across the validation workbook,
# plus collapsing a rich-text note into plain text and switching which sheet/
# cell is active when the file is saved.

$wb = $excel.ActiveWorkbook

$wsOverview  = $wb.Worksheets.Item("Environment Exposure Overview")
$wsBreakdown = $wb.Worksheets.Item("Environment Exposure Breakdown")

# ---------------------------------------------------------------------------
# 1. Breakdown sheet: column headers -- drop the redundant "BioGears " prefix
#    wording and rebrand to "Engine "
# ---------------------------------------------------------------------------
$wsBreakdown.Range("G1").Value  = "Engine HeartRate`n(BPM)"
$wsBreakdown.Range("K1").Value  = "Engine MeanArterialPressure`n(mmHg)"
$wsBreakdown.Range("M1").Value  = "Engine SystolicArterialPressure`n(mmHg)"
$wsBreakdown.Range("O1").Value  = "Engine DiastolicArterialPressure`n(mmHg)"
$wsBreakdown.Range("Q1").Value  = "Engine RespirationRate`n(Breaths/min)"
$wsBreakdown.Range("Y1").Value  = "Engine CoreTemperature`n(C)"
$wsBreakdown.Range("AA1").Value = "Engine SkinTemperature`n(C)"

# ---------------------------------------------------------------------------
# 2. Breakdown sheet: segment-0 note -- collapse the multi-run rich text into
#    plain text and rebrand "biogears documentation" -> "Engine documentation"
# ---------------------------------------------------------------------------
$wsBreakdown.Range("E2").Value = "Standard initialization buffer for scenarios. At the end of this segment this patient is in a resting physiological state. For validation references see the Engine documentation on resting physiology validation."

# ---------------------------------------------------------------------------
# 3. Overview sheet: scenario narrative -- "BioGears(R) physiology engine" -> "Pulse physiology engine"
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = "The environment exposure scenario simulates the physiology of an adult woman who is exposed to cold weather long enough to develop hypothermia. This scenario highlights the ability of the Pulse physiology engine to siimulate physiology when the body is exposed to an abnormal environment."

# ---------------------------------------------------------------------------
# 4. Switch the active sheet/selection: "Overview" becomes the active tab
#    with R10 selected (previously "Breakdown" was active with Y3 selected).
# ---------------------------------------------------------------------------
$wsOverview.Activate()
$wsOverview.Range("R10").Select()
